$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.444.47'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '3.774.58'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.59'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.15'
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.157'
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.39'
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000247'
$ws.Range('E12').Value = '  -2.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.43'
$ws.Range('E13').Value = '  -1.76%  '
$ws.Range('D14').Value = '4.409.91'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = '3.772.99'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '67.514.62'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.24'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.02'
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '459.42'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.69'
$ws.Range('E21').Value = '  -3.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.693'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000145'
$ws.Range('E23').Value = '  -5.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.22'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.95'
$ws.Range('E25').Value = '  -1.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.08'
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.92'
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').Value = '3.922.72'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.43'
$ws.Range('E30').Value = '  +3.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.60'
$ws.Range('E31').Value = '  -6.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.18'
$ws.Range('E32').Value = '  -2.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.99'
$ws.Range('E33').Value = '  -2.19%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.95'
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0982'
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.984'
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.21'
$ws.Range('E39').Value = '  -4.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.73'
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '47.38'
$ws.Range('E43').Value = '  -1.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.26'
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.81'
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.36'
$ws.Range('E48').Value = '  +8.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '26.95'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.83'
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '389.53'
$ws.Range('E51').Value = '  -0.05%  '

$ws.Range('D2:D51').ClearFormats()
